# Commit: "added an option to configure the cost function so that it would
# vary prelims/finals/prelim-finals"
#
# - Flip several 0/1 selection values on the "Strategy" sheet.
# - Add a new "Scoring" worksheet (same layout as "Strategy": column event
#   headers in row 1, swimmer names in column A) holding the numeric scoring
#   weights used by the (new) configurable cost function, and make it the
#   active sheet.

$wb = $excel.ActiveWorkbook
$strategy = $wb.Worksheets.Item("Strategy")

# ---------------------------------------------------------------------
# 1) Update the 0/1 selection values on the "Strategy" sheet
# ---------------------------------------------------------------------
$strategyUpdates = @{
    "M2" = 1; "O2" = 0;
    "H3" = 1; "P3" = 0;
    "C4" = 0; "M4" = 1; "N4" = 1; "Q4" = 0;
    "E5" = 1; "J5" = 1; "N5" = 0; "O5" = 0;
    "G6" = 1; "L6" = 1; "N6" = 0; "O6" = 0; "S6" = 0;
    "I7" = 1; "N7" = 0;
    "E8" = 0; "M8" = 1; "O8" = 0; "S8" = 1;
    "C9" = 1;
    "F10" = 1; "G10" = 1; "J10" = 1; "N10" = 0; "O10" = 0;
    "F11" = 1; "K11" = 1; "L11" = 0; "M11" = 0; "O11" = 1; "R11" = 0;
    "D12" = 1; "N12" = 0;
    "B13" = 1; "L13" = 1; "S13" = 0;
}

foreach ($addr in $strategyUpdates.Keys) {
    $strategy.Range($addr).Value = $strategyUpdates[$addr]
}

# ---------------------------------------------------------------------
# 2) Add the new "Scoring" worksheet right after "Strategy"
# ---------------------------------------------------------------------
$scoring = $wb.Worksheets.Add($null, $strategy)
$scoring.Name = "Scoring"

# Header row (row 1) - same column labels as Strategy
$headers = @("FR50m","FR100m","FR200m","FR400m","BR50m","BR100m","BA100m","BA50m","FLY50m","FLY100m","IM100m","IM200m","FRRelay4P50","FRRelay4P100","IMRelay4P50_FR","IMRelay4P50_BR","IMRelay4P50_BA","IMRelay4P50_FLY")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 2 # starts at column B
    $scoring.Cells.Item(1, $col).Value = $headers[$i]
}

# Row labels (column A) - same swimmer names as Strategy
$names = @("Miles Huang","Curtis Wong","King Wah","Justin Choi","Aaron Wu","Frank Zhou","Alan Wang","Alan Sun","Bernard Ip","Kan KikuchiYuan","Jerry Zheng","Aaron Sun")
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $scoring.Cells.Item($row, 1).Value = $names[$i]
}

# Scoring weight values for rows 2-13 (columns B-S)
$scoringData = @(
    @(0, 0, 4, 16, 0, 0, 10, 13, 0, 1, 0, 7, 0, 0, 0, 0, 3.25, 0),
    @(12, 4, 7, 11, 0, 0, 10, 5, 11, 9, 11, 7, 3, 1, 3, 0, 1.25, 2.75),
    @(12, 10, 1, 10, 16, 16, 10, 5, 7, 9, 5, 12, 3, 2.5, 3, 4, 1.25, 1.75),
    @(0, 0, 0, 11, 0, 0, 0, 0, 2, 0, 0, 0, 0, 0, 0, 0, 0, 0.5),
    @(0, 0, 0, 10, 0, 2, 0, 0, 1, 0, 5, 0, 0, 0, 0, 0, 0, 0.25),
    @(0, 0, 4, 16, 0, 0, 10, 5, 0, 0, 3, 12, 0, 0, 0, 0, 1.25, 0),
    @(1, 7, 13, 16, 3, 0, 0, 0, 10, 13, 7, 12, 0.25, 1.75, 0.25, 0.75, 0, 2.5),
    @(13, 2, 0, 2, 11, 0, 0, 0, 0, 0, 0, 0, 3.25, 0.5, 3.25, 2.75, 0, 0),
    @(0, 0, 0, 7, 9, 1, 0, 10, 7, 0, 0, 0, 0, 0, 0, 2.25, 2.5, 1.75),
    @(13, 16, 16, 16, 16, 16, 16, 16, 16, 16, 16, 16, 3.25, 4, 3.25, 4, 4, 4),
    @(0, 0, 9, 12, 0, 0, 0, 0, 0, 0, 0, 2, 0, 0, 0, 0, 0, 0),
    @(7, 0, 0, 7, 0, 0, 0, 0, 0, 9, 4, 0, 1.75, 0, 1.75, 0, 0, 0)
)

for ($r = 0; $r -lt $scoringData.Length; $r++) {
    $rowVals = $scoringData[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $scoring.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}

# Header + name cells share the same header style (s="1") already used on
# the Strategy sheet - copy that formatting over.
$strategy.Range("B1:S1").Copy()
$scoring.Range("B1:S1").PasteSpecial(-4122) # xlPasteFormats
$strategy.Range("A2:A13").Copy()
$scoring.Range("A2:A13").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Selection / active sheet bookkeeping
# ---------------------------------------------------------------------
$strategy.Range("S8").Select()
$scoring.Range("A1:S13").Select()
